$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Duplicate_Management")

# New DTC_Status column values (E2:E8), mirrors the Status column (C2:C8)
$ws.Range("E2").Value = "New"
$ws.Range("E3").Value = "In Progress"
$ws.Range("E4").Value = "Merged"
$ws.Range("E5").Value = "Suppressed"
$ws.Range("E6").Value = "Not a duplicate"
$ws.Range("E7").Value = "Unable to process"
$ws.Range("E8").Value = "Deleted"

# New header cells E1:G1
$ws.Range("E1").Value = "DTC_Status"
$ws.Range("F1").Value = "External Notes"
$ws.Range("G1").Value = "Internal Notes"

# Internal Notes column (G2:G8)
$ws.Range("G2:G8").Value = "Internal Notes"

# External Notes / Testing column (F2:F8)
$ws.Range("F2").Value = "Testing 123"
$ws.Range("F3").Value = "Testing 124"
$ws.Range("F4").Value = "Testing 125"
$ws.Range("F5").Value = "Testing 126"
$ws.Range("F6").Value = "Testing 127"
$ws.Range("F7").Value = "Testing 128"
$ws.Range("F8").Value = "Testing 129"

# Phone number (A2) and Location Number (D2:D8) content update
$ws.Range("A2").Value = "+1 610-488-2411"
$ws.Range("D2:D8").Value = "'9000223167"

# Update the view: scroll so column B is the leftmost visible, select H1:H1048576
$ws.Range("H1:H1048576").Select()

# Finally, activate the Duplicate_Management sheet/tab
$ws.Activate()
